# Fruta / hortaliza, semanal
# A new weekly price observation was added to the dataset. This inserts a new
# row at position 43 (shifting all subsequent rows down by one) and fills it
# in with the new record's data, mirroring the structure of the surrounding
# rows in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 43; everything below (old rows
# 43-119) shifts down to 44-120, and the sheet's used range grows to R120.
$ws.Rows("43:43").Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A43").Value = 10
$ws.Range("B43").Value = "Vega Modelo de Temuco"
$ws.Range("C43").Value = "La Araucanía"
$ws.Range("D43").Value = 45246
$ws.Range("E43").Value = 9
$ws.Range("F43").Value = 300000000
$ws.Range("G43").Value = "Espárragos"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 650
$ws.Range("K43").Value = 1800
$ws.Range("L43").Value = 1800
$ws.Range("M43").Value = 1800
$ws.Range("N43").Value = "`$/kilo"
$ws.Range("O43").Value = "Región del Maule"
$ws.Range("P43").Value = 1800
$ws.Range("Q43").Value = 1
$ws.Range("R43").Value = "Hortaliza"

# Make sure the date cell keeps the same date display format as the rest of
# column D (style index 2 in the original file, numFmt "YYYY-MM-DD HH:MM:SS").
$ws.Range("D43").NumberFormat = $ws.Range("D44").NumberFormat()
